$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.911.66"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.73%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.888.09"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.99%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.016"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.97%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "336.52"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.86%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.016"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.87%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4796"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.12%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.43"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.16%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08024"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.00%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.018"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.65%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.73"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.03%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.885.29"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.09%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.982"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.77%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.215"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.32%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.019"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.20%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.06808"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.19%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "88.20"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.64%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.00001048"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.00%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.01"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.62%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.016"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.90%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "27.882.62"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.59%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.482"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.97"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.35%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.358"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.96%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.105.11"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.66%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "159.99"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.39%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.99"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.90%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.090"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.20%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.483"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.85%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "121.63"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.62%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09617"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.93%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9588"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.95%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.650"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.10%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.331"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.24%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.355"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.87%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06085"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.05%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02231"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.57%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.200"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.94%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.144"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.18%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5924"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.42%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1902"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.04%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.35"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.49%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.273"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5665"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.90%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.13"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.23%  "
$ws.Range("B47").Value = "PancakeSwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.382"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.71%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.932"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.40%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06855"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.33%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "113.67"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.41%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.070"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.15%  "
